# Replace the synthetic/placeholder OHLC + shares_outstanding + fixed_ticker data
# (rows 2-43 of Sheet1) with the corrected values, and normalize every
# "fixed_ticker" entry (column I) to "SSTK" instead of the stray unrelated
# ticker codes that had leaked in from other companies' data pulls.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{ Row = 2; D = 56.68370742830358; E = 55.15546798706055; F = 60.72087078787232; G = 54.27285445170207; H = 35474129 },
    @{ Row = 3; D = 48.47863105179211; E = 43.66509246826172; F = 48.67476983941955; G = 39.11307018988773; H = 35474129 },
    @{ Row = 4; D = 24.7541709597851; E = 23.27496719360352; F = 28.11302387092217; G = 22.4740733185362; H = 35474129 },
    @{ Row = 5; D = 25.83293026294953; E = 23.61003875732422; F = 27.11599549060297; G = 22.63752528237314; H = 35474129 },
    @{ Row = 6; D = 29.59223320340016; E = 33.52315521240234; F = 35.38645974944534; G = 29.21630388759681; H = 35474129 },
    @{ Row = 7; D = 37.26610459661102; E = 45.02170562744141; F = 47.23642301160957; G = 36.78393336560024; H = 35474129 },
    @{ Row = 8; D = 52.13168587628555; E = 48.20893859863281; F = 52.96526706186381; G = 47.71042222072429; H = 35474129 },
    @{ Row = 9; D = 39.16210769653555; E = 43.96747589111328; F = 44.23716636986193; G = 38.55735233043603; H = 35474129 },
    @{ Row = 10; D = 33.91543421411985; E = 35.32925796508789; F = 35.9013261536495; G = 33.56402103189718; H = 35474129 },
    @{ Row = 11; D = 36.28541570563405; E = 34.43845367431641; F = 37.38868719554388; G = 34.20145359744049; H = 35474129 },
    @{ Row = 12; D = 27.22223797155721; E = 31.86415672302246; F = 34.90428432629136; G = 26.88717003536962; H = 35474129 },
    @{ Row = 13; D = 35.25569450915763; E = 36.17100143432617; F = 39.52985465641906; G = 35.21483317739042; H = 35474129 },
    @{ Row = 14; D = 39.13757829452138; E = 34.43845367431641; F = 42.61901845460596; G = 34.16059226772524; H = 35474129 },
    @{ Row = 15; D = 38.73713708914858; E = 37.65020751953125; F = 42.18588419399941; G = 37.29062341461182; H = 35474129 },
    @{ Row = 16; D = 47.65695943814214; E = 35.422119140625; F = 48.12486492410298; G = 32.1467840442968; H = 35474129 },
    @{ Row = 17; D = 30.6304226965719; E = 34.66826629638672; F = 36.12396988052056; G = 29.37401473757547; H = 35474129 },
    @{ Row = 18; D = 40.72503587182806; E = 35.04952621459961; F = 41.88613277145352; G = 32.8399770571012; H = 35474129 },
    @{ Row = 19; D = 34.44298406864268; E = 33.24722671508789; F = 34.93688319840273; G = 31.97348841984246; H = 35474129 },
    @{ Row = 20; D = 31.48825087268828; E = 35.16217041015625; F = 35.90734913160824; G = 29.40000983031588; H = 35474129 },
    @{ Row = 21; D = 37.35439215994818; E = 37.54502105712891; F = 39.26933466298065; G = 36.44457813272712; H = 35474129 },
    @{ Row = 22; D = 27.31375852517522; E = 33.07593536376953; F = 34.6949158891527; G = 25.33790785861263; H = 35474129 },
    @{ Row = 23; D = 30.57062351176283; E = 47.51738357543945; F = 47.93711714740178; G = 30.37824770945463; H = 35474129 },
    @{ Row = 24; D = 45.87541563538974; E = 57.4210319519043; F = 62.67622389770459; G = 45.10336717387623; H = 35474129 },
    @{ Row = 25; D = 63.97114315621747; E = 57.15541076660156; F = 67.77916041131625; G = 57.14662110415512; H = 35474129 },
    @{ Row = 26; D = 78.87658762310899; E = 76.85772705078125; F = 92.18871835411004; G = 76.12599897141048; H = 35474129 },
    @{ Row = 27; D = 86.9186214789853; E = 95.86053466796876; F = 96.27582242903752; G = 85.16027513993899; H = 35474129 },
    @{ Row = 28; D = 100.5588904896665; E = 107.2421646118164; F = 113.6244667078085; G = 98.5406321213258; H = 35474129 },
    @{ Row = 29; D = 98.5805734042838; E = 85.99638366699219; F = 99.0239908663052; G = 74.24582092342513; H = 35474129 },
    @{ Row = 30; D = 83.49962549937156; E = 67.33324432373047; F = 84.80680919924085; G = 65.27910238642835; H = 35474129 },
    @{ Row = 31; D = 50.98776527883513; E = 50.44315719604492; F = 56.34456073328239; G = 44.67567490747216; H = 35474129 },
    @{ Row = 32; D = 45.4070422242993; E = 44.86007690429688; F = 48.41983233936182; G = 40.99545914188635; H = 35474129 },
    @{ Row = 33; D = 48.73925721882901; E = 67.81153106689453; F = 70.84760196017547; G = 47.44194732157453; H = 35474129 },
    @{ Row = 34; D = 65.15347619411369; E = 60.57844924926758; F = 68.67968581362894; G = 58.77918082245395; H = 35474129 },
    @{ Row = 35; D = 44.55843178205527; E = 46.78635406494141; F = 53.12456128813761; G = 43.85822743778151; H = 35474129 },
    @{ Row = 36; D = 34.73243046815878; E = 37.2309684753418; F = 38.4573574395426; G = 30.93428457770854; H = 35474129 },
    @{ Row = 37; D = 43.8367249439768; E = 43.2474479675293; F = 47.5197095591415; G = 41.70059414817069; H = 35474129 },
    @{ Row = 38; D = 42.6894199636056; E = 39.57597351074219; F = 43.05080150995978; G = 37.15749265772229; H = 35474129 },
    @{ Row = 39; D = 36.03350380922711; E = 41.30123138427734; F = 43.05713820069029; G = 32.12006519420992; H = 35474129 },
    @{ Row = 40; D = 33.42040529528817; E = 30.23571395874023; F = 33.42040529528817; G = 27.277156061433; H = 35474129 },
    @{ Row = 41; D = 29.17324605295624; E = 28.0793685913086; F = 36.90648587484644; G = 25.96770762983717; H = 35474129 },
    @{ Row = 42; D = 17.82369365106175; E = 15.42658042907715; F = 18.53896094174895; G = 13.87039063364142; H = 35474129 },
    @{ Row = 43; D = 18.59357894412627; E = 18.85934257507324; F = 20.9559198270967; G = 17.61911417473937; H = 35474129 }
)

foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = "SSTK"
}
